# Iterative (fixed-point) root-finding table: renamed/updated with descriptive
# values for section 1 and 2, using a relaxation factor w with 0 < w < 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0", "3.0",                 "-17.6902860206768",    "1.000005"),
    @("1", "20.6902860206768",    "-13.5173641896727",    "17.6902860206768"),
    @("2", "34.2076502103495",    "-3.51949495555342",    "13.5173641896727"),
    @("3", "37.7271451659029",    "-0.685513898070205",   "3.51949495555342"),
    @("4", "38.4126590639731",    "-0.126050414764144",   "0.685513898070205"),
    @("5", "38.5387094787372",    "-0.0229327615412593",  "0.126050414764144"),
    @("6", "38.5616422402785",    "-0.0041641661015425",  "0.0229327615412629"),
    @("7", "38.56580640638",      "-0.0007558700112753",  "0.0041641661015461"),
    @("8", "38.5665622763913",    "-0.000137195062706",   "0.0007558700112753"),
    @("9", "38.566699471454",     "-2.49014605948616e-05","0.000137195062706"),
    @("10","38.5667243729146",    "-4.51970664983037e-06","2.49014605913089e-05"),
    @("11","38.5667288926213",    "-8.20343057483797e-07","4.51970664983037e-06")
)

# Force the target range to be stored as text so the literal numeric strings
# (including special forms like "3.0") are preserved exactly, rather than
# being reinterpreted as numbers by Excel.
$rng = $ws.Range("A2:D13")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
